# Auto-generated market-data refresh for the Leve profit sheets.
# Mirrors a scheduled runner pulling fresh Universalis prices into
# currentAveragePrice / LevePrice / LeveProfit columns (H-N) per sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 494.06668
$ws.Range("I28").Value = 494.06668
$ws.Range("K28").Value = 494.06668
$ws.Range("M28").Value = -9.066680000000019
$ws.Range("H107").Value = 1024.6923
$ws.Range("I107").Value = 1101.65
$ws.Range("J107").Value = 768.1667
$ws.Range("K107").Value = 1101.65
$ws.Range("L107").Value = 768.1667
$ws.Range("M107").Value = 818.3499999999999
$ws.Range("N107").Value = -4608.1667
$ws.Range("H125").Value = 4500.2354
$ws.Range("I125").Value = 3247
$ws.Range("J125").Value = 4885.846
$ws.Range("K125").Value = 29223
$ws.Range("L125").Value = 43972.61399999999
$ws.Range("M125").Value = -26763
$ws.Range("N125").Value = -48892.61399999999
$ws.Range("H132").Value = 4449.6
$ws.Range("I132").Value = 3927.48
$ws.Range("J132").Value = 7060.2
$ws.Range("K132").Value = 11782.44
$ws.Range("L132").Value = 21180.6
$ws.Range("M132").Value = -9252.440000000001
$ws.Range("N132").Value = -26240.6
$ws.Range("H136").Value = 99999.5
$ws.Range("J136").Value = 99999.5
$ws.Range("L136").Value = 99999.5
$ws.Range("N136").Value = -110199.5
$ws.Range("H141").Value = 1176.1072
$ws.Range("I141").Value = 1208.0416
$ws.Range("J141").Value = 984.5
$ws.Range("K141").Value = 3624.1248
$ws.Range("L141").Value = 2953.5
$ws.Range("M141").Value = 1555.8752
$ws.Range("N141").Value = -13313.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 112819.664
$ws.Range("I74").Value = 112819.664
$ws.Range("K74").Value = 112819.664
$ws.Range("M74").Value = -111945.664
$ws.Range("H77").Value = 112819.664
$ws.Range("I77").Value = 112819.664
$ws.Range("K77").Value = 564098.3200000001
$ws.Range("M77").Value = -559730.3200000001
$ws.Range("H97").Value = 5746.577
$ws.Range("I97").Value = 6166.3
$ws.Range("K97").Value = 6166.3
$ws.Range("M97").Value = -5670.3
$ws.Range("H132").Value = 56206.105
$ws.Range("I132").Value = 62442.234
$ws.Range("K132").Value = 187326.702
$ws.Range("M132").Value = -184796.702

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 59361.883
$ws.Range("I22").Value = 62947
$ws.Range("K22").Value = 62947
$ws.Range("M22").Value = -62774
$ws.Range("H51").Value = 59999
$ws.Range("J51").Value = 59999
$ws.Range("L51").Value = 59999
$ws.Range("N51").Value = -60981
$ws.Range("H94").Value = 3022.0908
$ws.Range("I94").Value = 2974.9333
$ws.Range("J94").Value = 3123.1428
$ws.Range("K94").Value = 2974.9333
$ws.Range("L94").Value = 3123.1428
$ws.Range("M94").Value = -2523.9333
$ws.Range("N94").Value = -4025.1428
$ws.Range("H99").Value = 74343.14
$ws.Range("I99").Value = 102880.5
$ws.Range("K99").Value = 102880.5
$ws.Range("M99").Value = -101382.5
$ws.Range("H105").Value = 3496.0303
$ws.Range("I105").Value = 3239.2
$ws.Range("J105").Value = 4298.625
$ws.Range("K105").Value = 3239.2
$ws.Range("L105").Value = 4298.625
$ws.Range("M105").Value = -1492.2
$ws.Range("N105").Value = -7792.625
$ws.Range("H124").Value = 129999
$ws.Range("J124").Value = 129999
$ws.Range("L124").Value = 129999
$ws.Range("N124").Value = -139819
$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1180.2941
$ws.Range("I16").Value = 769.5454999999999
$ws.Range("K16").Value = 769.5454999999999
$ws.Range("M16").Value = -482.5454999999999
$ws.Range("H31").Value = 2785.2173
$ws.Range("I31").Value = 2003.6842
$ws.Range("J31").Value = 6497.5
$ws.Range("K31").Value = 2003.6842
$ws.Range("L31").Value = 6497.5
$ws.Range("M31").Value = -1708.6842
$ws.Range("N31").Value = -7087.5
$ws.Range("H34").Value = 2785.2173
$ws.Range("I34").Value = 2003.6842
$ws.Range("J34").Value = 6497.5
$ws.Range("K34").Value = 2003.6842
$ws.Range("L34").Value = 6497.5
$ws.Range("M34").Value = -1801.6842
$ws.Range("N34").Value = -6901.5
$ws.Range("H113").Value = 1180.2941
$ws.Range("I113").Value = 769.5454999999999
$ws.Range("K113").Value = 769.5454999999999
$ws.Range("M113").Value = 1400.4545
$ws.Range("H132").Value = 2360.389
$ws.Range("I132").Value = 2149.6365
$ws.Range("K132").Value = 6448.9095
$ws.Range("M132").Value = -3918.9095
$ws.Range("H134").Value = 126930.375
$ws.Range("I134").Value = 126930.375
$ws.Range("K134").Value = 380791.125
$ws.Range("M134").Value = -378256.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 4036.6667
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H68").Value = 5899.8335
$ws.Range("I68").Value = 150
$ws.Range("K68").Value = 450
$ws.Range("M68").Value = 361
$ws.Range("H71").Value = 5899.8335
$ws.Range("I71").Value = 150
$ws.Range("K71").Value = 1350
$ws.Range("M71").Value = 2706

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4735.353
$ws.Range("I80").Value = 4923.091
$ws.Range("J80").Value = 4391.1665
$ws.Range("K80").Value = 4923.091
$ws.Range("L80").Value = 4391.1665
$ws.Range("M80").Value = -3925.091
$ws.Range("N80").Value = -6387.1665
$ws.Range("H83").Value = 4735.353
$ws.Range("I83").Value = 4923.091
$ws.Range("J83").Value = 4391.1665
$ws.Range("K83").Value = 24615.455
$ws.Range("L83").Value = 21955.8325
$ws.Range("M83").Value = -19623.455
$ws.Range("N83").Value = -31939.8325
$ws.Range("H102").Value = 3579.6667
$ws.Range("I102").Value = 3572.3635
$ws.Range("K102").Value = 3572.3635
$ws.Range("M102").Value = -1950.3635
$ws.Range("H122").Value = 4236.364
$ws.Range("I122").Value = 2074.125
$ws.Range("K122").Value = 6222.375
$ws.Range("M122").Value = -3772.375
$ws.Range("H132").Value = 41963.242
$ws.Range("I132").Value = 68476.7
$ws.Range("J132").Value = 4402.5
$ws.Range("K132").Value = 205430.1
$ws.Range("L132").Value = 13207.5
$ws.Range("M132").Value = -202900.1
$ws.Range("N132").Value = -18267.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9227.214
$ws.Range("I40").Value = 7082.5
$ws.Range("J40").Value = 10835.75
$ws.Range("K40").Value = 7082.5
$ws.Range("L40").Value = 10835.75
$ws.Range("M40").Value = -6946.5
$ws.Range("N40").Value = -11107.75
$ws.Range("H132").Value = 90799.78999999999
$ws.Range("I132").Value = 123569.7
$ws.Range("K132").Value = 370709.1
$ws.Range("M132").Value = -368179.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 58999.5
$ws.Range("J110").Value = 58999.5
$ws.Range("L110").Value = 58999.5
$ws.Range("N110").Value = -67179.5
$ws.Range("H132").Value = 25878.717
$ws.Range("I132").Value = 29778.36
$ws.Range("K132").Value = 89335.08
$ws.Range("M132").Value = -86805.08
